$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1027.384161038
$ws.Range("C1").Value = 9.8443438331496207
$ws.Range("G1").Value = 9.8418471971860697
$ws.Range("Q1").Value = 9.8584060415702996
$ws.Range("T1").Value = 0.202506901722942
$ws.Range("U1").Value = 9.9381617420716495
$ws.Range("AE1").Value = 9.8584060415702996
$ws.Range("AH1").Value = 8.2233768111590599
$ws.Range("AI1").Value = 9.9381617420716495
$ws.Range("AL1").Value = 5.3002636771839402
$ws.Range("AP1").Value = 5.3002636771839402
$ws.Range("B2").Value = 2054.768322076
$ws.Range("C2").Value = 19.762205894692976
$ws.Range("G2").Value = 19.743004418991941
$ws.Range("Q2").Value = 28.467171455087001
$ws.Range("R2").Value = 8.2202576742737694
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 0.36298914529252702
$ws.Range("U2").Value = 28.460717719406901
$ws.Range("V2").Value = 3
$ws.Range("AE2").Value = 22.7312018312996
$ws.Range("AH2").Value = 22.413708315094699
$ws.Range("AI2").Value = 22.7457347018306
$ws.Range("AL2").Value = 9.9033685267320362
$ws.Range("AP2").Value = 9.9033685267320362
$ws.Range("B3").Value = 3082.152483114
$ws.Range("C3").Value = 29.115972251621969
$ws.Range("G3").Value = 29.085727412886488
$ws.Range("Q3").Value = 52.018620258839903
$ws.Range("R3").Value = 21.937090365042
$ws.Range("S3").Value = 8
$ws.Range("T3").Value = 0.31410029503186598
$ws.Range("U3").Value = 51.833036634750897
$ws.Range("V3").Value = 8
$ws.Range("AE3").Value = 35.853705045367001
$ws.Range("AH3").Value = 66.001699949274098
$ws.Range("AI3").Value = 35.7370575706423
$ws.Range("AL3").Value = 12.407563508323724
$ws.Range("AP3").Value = 12.407563508323724
$ws.Range("B4").Value = 4109.536644152
$ws.Range("C4").Value = 41.924944702572205
$ws.Range("G4").Value = 41.832294235166387
$ws.Range("Q4").Value = 76.572763806806506
$ws.Range("R4").Value = 35.6616051636694
$ws.Range("S4").Value = 13
$ws.Range("T4").Value = 0.67058119855312703
$ws.Range("U4").Value = 76.510332493623594
$ws.Range("V4").Value = 13
$ws.Range("AE4").Value = 49.649240883804403
$ws.Range("AH4").Value = 87.610868060548697
$ws.Range("AI4").Value = 49.593963830450498
$ws.Range("AL4").Value = 16.860545565823504
$ws.Range("AP4").Value = 16.860545565823504
$ws.Range("B5").Value = 5136.9208051900005
$ws.Range("C5").Value = 52.814924323620332
$ws.Range("G5").Value = 52.654640661424608
$ws.Range("Q5").Value = 109.45608293239199
$ws.Range("R5").Value = 57.6494277609414
$ws.Range("S5").Value = 21
$ws.Range("T5").Value = 0.53202091297113097
$ws.Range("U5").Value = 109.24047552071799
$ws.Range("V5").Value = 21
$ws.Range("AE5").Value = 75.875051475845297
$ws.Range("AH5").Value = 97.878079306209003
$ws.Range("AI5").Value = 75.655304043969196
$ws.Range("AL5").Value = 20.319508511507408
$ws.Range("AP5").Value = 20.319508511507408
$ws.Range("B6").Value = 6164.304966228
$ws.Range("C6").Value = 66.441826425685477
$ws.Range("G6").Value = 65.586733890776046
$ws.Range("Q6").Value = 138.80967346852
$ws.Range("R6").Value = 76.881476792806197
$ws.Range("S6").Value = 28
$ws.Range("T6").Value = 0.93956073911250204
$ws.Range("U6").Value = 138.66836843029
$ws.Range("V6").Value = 28
$ws.Range("AE6").Value = 96.270181883432699
$ws.Range("AH6").Value = 129.88522636229001
$ws.Range("AI6").Value = 96.300130319730897
$ws.Range("AL6").Value = 23.661184528677516
$ws.Range("AP6").Value = 23.661184528677516
$ws.Range("B7").Value = 7191.6891272659996
$ws.Range("C7").Value = 76.751976545525295
$ws.Range("G7").Value = 75.011646641741478
$ws.Range("Q7").Value = 170.745448490825
$ws.Range("R7").Value = 98.935511330031503
$ws.Range("S7").Value = 36
$ws.Range("T7").Value = 0.99355299274825304
$ws.Range("U7").Value = 170.19674237878701
$ws.Range("V7").Value = 36
$ws.Range("AE7").Value = 127.131878535421
$ws.Range("AH7").Value = 80.988331466009697
$ws.Range("AI7").Value = 127.220184130874
$ws.Range("AL7").Value = 27.182327681147065
$ws.Range("AP7").Value = 27.182327681147065
$ws.Range("B8").Value = 8219.073288304
$ws.Range("C8").Value = 90.488391899777923
$ws.Range("G8").Value = 89.394662569256539
$ws.Range("Q8").Value = 204.19993319915301
$ws.Range("R8").Value = 120.96653357353399
$ws.Range("S8").Value = 44
$ws.Range("T8").Value = 0.76531969159925595
$ws.Range("U8").Value = 203.793290075279
$ws.Range("V8").Value = 44
$ws.Range("AE8").Value = 139.50809804821199
$ws.Range("AH8").Value = 201.921388704088
$ws.Range("AI8").Value = 139.30322060703801
$ws.Range("AL8").Value = 34.878759329678999
$ws.Range("AP8").Value = 34.878759329678999
$ws.Range("B9").Value = 9246.4574493420005
$ws.Range("C9").Value = 97.767615081483342
$ws.Range("G9").Value = 95.817868587472077
$ws.Range("Q9").Value = 254.61607064427801
$ws.Range("R9").Value = 148.564946060436
$ws.Range("S9").Value = 54
$ws.Range("T9").Value = 0.88202359938876096
$ws.Range("U9").Value = 253.68559406050599
$ws.Range("V9").Value = 54
$ws.Range("AE9").Value = 175.54783964201047
$ws.Range("AH9").Value = 143.90268959960247
$ws.Range("AI9").Value = 178.32077633946301
$ws.Range("AL9").Value = 42.033180460677798
$ws.Range("AP9").Value = 42.033180460677826
$ws.Range("B10").Value = 10273.841610380001
$ws.Range("C10").Value = 109.5009830679254
$ws.Range("G10").Value = 108.29160753283814
$ws.Range("Q10").Value = 546.52197946495698
$ws.Range("R10").Value = 176.22350135519801
$ws.Range("S10").Value = 64
$ws.Range("T10").Value = 0.807347115182526
$ws.Range("U10").Value = 545.95170437980005
$ws.Range("V10").Value = 64
$ws.Range("AE10").Value = 191.76053535459999
$ws.Range("AH10").Value = 225.540570042019
$ws.Range("AI10").Value = 191.77803619477999
$ws.Range("AL10").Value = 48.94572719134429
$ws.Range("AP10").Value = 48.94572719134429

$ws.Range("AE13").Select()
